$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns remain text so values are not
# reinterpreted as numbers (losing trailing zeros / formatting like "34.343.45").
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D2").Value = "34.343.45"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.804.85"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "227.50"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").Value = "0.575"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "36.15"
$ws.Range("E8").Value = "  +11.01%  "
$ws.Range("E9").Value = "  +2.18%  "
$ws.Range("D10").Value = "0.0694"
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("D11").Value = "0.0965"
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("D12").Value = "2.065.34"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("E13").Value = "  +6.25%  "
$ws.Range("D14").Value = "1.819.29"
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("D15").Value = "0.645"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").Value = "4.51"
$ws.Range("E16").Value = "  +5.33%  "
$ws.Range("D17").Value = "34.341.48"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "69.09"
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("D19").Value = "245.49"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").Value = "0.0₃0795"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "11.59"
$ws.Range("E21").Value = "  +3.70%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").Value = "171.75"
$ws.Range("E24").Value = "  +3.64%  "
$ws.Range("D25").Value = "2.12"
$ws.Range("E25").Value = "  +2.72%  "
$ws.Range("D26").Value = "8.00"
$ws.Range("E26").Value = "  +9.70%  "
$ws.Range("E28").Value = "  +2.65%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  +0.74%  "
$ws.Range("E31").Value = "  +1.42%  "
$ws.Range("D32").Value = "3.84"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("D35").Value = "1.390.91"
$ws.Range("E35").Value = "  -1.29%  "
$ws.Range("D36").Value = "0.673"
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("D37").Value = "2.47"
$ws.Range("E37").Value = "  -5.73%  "
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("E40").Value = "  +10.60%  "
$ws.Range("D41").Value = "0.966"
$ws.Range("E41").Value = "  +3.00%  "
$ws.Range("D42").Value = "82.42"
$ws.Range("E42").Value = "  -1.96%  "
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("D44").Value = "2.42"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").Value = "13.47"
$ws.Range("E45").Value = "  -3.01%  "
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").Value = "0.0503"
$ws.Range("E47").Value = "  -4.40%  "
$ws.Range("D48").Value = "1.965.93"
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("D49").Value = "104.50"
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").Value = "0.0₆0129"
$ws.Range("E51").Value = "  +0.18%  "
